# [Kadastro App] Yeni kayit eklendi: 3000
#
# Appends a new record row (row 60) to both the "Kayitlar" sheet and the
# "Erdemli" sheet, which carry the same tabular data (Kayit No, Tarih,
# Birim, Parsel Sayisi, Is, Personeller). All the columns in this table
# are stored as text, even the numeric-looking ones (Kayit No, Parsel
# Sayisi), so we force a Text number format before writing the values to
# avoid Excel auto-coercing them into numbers/dates, then restore the
# default "Normal" style so no stray formatting is left behind.

$wb = $excel.ActiveWorkbook

$newRow = @("3000", "2025-09-11", "Erdemli", "2", "MAHKEME KARARI", "EMİNE ALANLI KIRCILI (K.Mühendisi), AYHAN KARADAYI (K.Teknisyeni)")

$sheetNames = @("Kayitlar", "Erdemli")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    $lastRow = $ws.Cells.Item(1, 1).EntireColumn.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
    $targetRow = $lastRow + 1

    $rng = $ws.Range($ws.Cells.Item($targetRow, 1), $ws.Cells.Item($targetRow, 6))
    $rng.NumberFormat = "@"

    for ($col = 1; $col -le 6; $col++) {
        $ws.Cells.Item($targetRow, $col).Value = $newRow[$col - 1]
    }

    $rng.Style = "Normal"
}

Write-Host "Added record 3000 to Kayitlar and Erdemli sheets"
